$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1: "Save", formatted like the existing header cells (e.g. G1)
$ws.Cells.Item(1, 8).Value = "Save"
$ws.Cells.Item(1, 7).Copy()
$ws.Cells.Item(1, 8).PasteSpecial(-4122)   # xlPasteFormats

# New data column H2:H8, filled with 0 (plain/unstyled numbers like the rest of the data rows)
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}

$wb.Save()
